$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.085.21'
$ws.Range("E2").Value = '  -1.92%  '
$ws.Range("D3").Value = '1.995.81'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").Value = '''1.013'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''330.25'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '''1.011'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = '''0.4968'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").Value = '''55.06'
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").Value = '''0.08905'
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("E11").Value = '  -3.47%  '
$ws.Range("D12").Value = '''22.91'
$ws.Range("E12").Value = '  -3.09%  '
$ws.Range("D13").Value = '2.005.74'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '''8.000'
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '''6.414'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '''92.47'
$ws.Range("E17").Value = '  -4.00%  '
$ws.Range("D18").Value = '''0.00001106'
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = '''0.06721'
$ws.Range("E19").Value = '  +0.99%  '
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '''5.974'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = '29.123.63'
$ws.Range("E23").Value = '  -1.90%  '
$ws.Range("D24").Value = '''11.96'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '''2.294'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").Value = '2.246.12'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").Value = '''20.81'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = '''156.85'
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").Value = '''6.263'
$ws.Range("E29").Value = '  -4.55%  '
$ws.Range("D30").Value = '''2.247'
$ws.Range("E30").Value = '  -4.68%  '
$ws.Range("D31").Value = '''127.00'
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("E32").Value = '  -2.17%  '
$ws.Range("D33").Value = '''0.09879'
$ws.Range("E33").Value = '  -1.18%  '
$ws.Range("D34").Value = '''1.532'
$ws.Range("E34").Value = '  -4.48%  '
$ws.Range("D35").Value = '''5.826'
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").Value = '''3.745'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").Value = '''0.02417'
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("D38").Value = '''1.310'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").Value = '''9.070'
$ws.Range("E39").Value = '  -6.45%  '
$ws.Range("D40").Value = '''0.06362'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("D41").Value = '''0.6467'
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("E43").Value = '  -5.40%  '
$ws.Range("D44").Value = '''1.011'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("E45").Value = '  +4.95%  '
$ws.Range("D46").Value = '''0.6177'
$ws.Range("E46").Value = '  -3.12%  '
$ws.Range("D47").Value = '''13.36'
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("D48").Value = '''2.167'
$ws.Range("E48").Value = '  -2.78%  '
$ws.Range("D49").Value = '''0.00000000356'
$ws.Range("E49").Value = '  +10.41%  '
$ws.Range("D50").Value = '''3.493'
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = '''2.183'
$ws.Range("E51").Value = '  +6.53%  '
